# Update CDC data: add WHO data reports 1 through 40 (column M, "26 mar data")
# and fix a handful of existing column B values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header label for column M
$ws.Range("M1").Value = "26 mar data"

# New column M values for rows 2-76 (one value per row, in order)
$mValues = @(0,0,2,0,1,0,0,0,1,1,2,0,1,3,0,0,3,1,2,0,1,4,0,0,0,2,1,1,0,4,6,0,4,3,7,3,13,14,7,12,21,20,23,48,40,74,71,102,97,173,188,216,247,274,365,385,548,782,767,993,1341,1641,1634,1905,1880,1246,977,896,457,179,90,59,7,1,1)

for ($i = 0; $i -lt $mValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 13).Value = $mValues[$i]
}

# Corrections to existing column B ("cases") values
$ws.Range("B37").Value = 5
$ws.Range("B49").Value = 103
$ws.Range("B51").Value = 273
$ws.Range("B68").Value = 997
